$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (the current "modifiedDate" column),
# shifting existing F:P columns to G:Q.
$ws.Columns("F:F").Insert()

# Give the newly inserted column the same (approximate) width as column E,
# since in the target file it shares the same customWidth value.
$ws.Columns("F:F").ColumnWidth = 14.022135416666666

# New header label for the inserted column.
$ws.Range("F1").Value = "ongkosKirimBeli"

# Update the active cell/selection to reflect where the user ended up
# after performing the insert (on the newly created column, row 2).
$ws.Range("F2").Select()
